$d = $word.ActiveDocument

$importStyle = $d.Styles("ImportTok")
$importStyle.Font.Color = 32768
$importStyle.Font.Bold = $true

$builtInStyle = $d.Styles("BuiltInTok")
$builtInStyle.Font.Color = 32768
